$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Light gray fill (D9D9D9) for the blank separator row (row 6)
$grayColor = 217 + (217 * 256) + (217 * 65536)
$ws.Range("A6:E6").Interior.Color = $grayColor

# Row 7 - Daniel Williams
$ws.Cells.Item(7,1).Value = "daniel.williams@example.com"
$ws.Cells.Item(7,2).Value = "416 715 6897"
$ws.Cells.Item(7,3).Value = "Daniel Williams"
$ws.Cells.Item(7,4).Value = "766 Birch Blvd., L8K7J6, Charletown, PE, US"
$ws.Cells.Item(7,5).Value = "Unit ID : unit-4970    Unit Type : Apartment    Owner : Great   "

# Rows 8-16 - Jane Davis (repeated)
for ($r = 8; $r -le 16; $r++) {
    $ws.Cells.Item($r,1).Value = "jane.davis@example.com"
    $ws.Cells.Item($r,2).Value = "416 715 6897"
    $ws.Cells.Item($r,3).Value = "Jane Davis"
    $ws.Cells.Item($r,4).Value = "4947 Maple St., L5K7J6, Regina, SK, US"
    $ws.Cells.Item($r,5).Value = "Unit ID : unit-4493    Unit Type : Apartment    Owner : Average   "
}
